$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.74%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.33%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.567"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.42%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08043"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.65%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.906"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.06%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.288"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.04%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9458"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.03%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1167"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.79%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1842"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.66%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09713"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.18%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04397"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.96%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1068"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.08%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001282"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.61%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005957"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.22%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.416"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.44%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3493"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.24%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.581"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "9.80%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1379"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.71%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2509"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.59%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04215"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.10%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001249"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.67%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004335"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.22%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001262"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.14%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003997"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.26%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02640"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-6.15%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05490"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.15%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007595"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.92%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1396"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.47%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008055"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-17.87%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002012"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.25%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008619"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.96%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006919"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.46%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.26%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.002275"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.26%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.005248"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "52.92%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.26%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.26%"

